# occurences_clients Z.xlsx - "hardware clients et instance clients"
# Rename client/company names in the COMPANY_BREAKDOWN column (J):
#   ETHIAS,AME LIFE              -> ETHIAS,AME-LIFE
#   RESA,PLG                     -> RESA,PROVINCE LIEGE
#   ETHIAS,AME LIFE,RESA,PLG     -> ETHIAS,AME-LIFE,RESA,PROVINCE LIEGE

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the company names wherever they occur in the sheet.
$ws.Cells.Replace("ETHIAS,AME LIFE", "ETHIAS,AME-LIFE")
$ws.Cells.Replace("RESA,PLG", "RESA,PROVINCE LIEGE")

# Row 14 (IDA1 / "IMS développement mutualisé") gets its own, freshly typed
# value for the combined company breakdown (no trailing space), distinct
# from row 13 (IHA1) which keeps the shared/previously-edited text.
$ws.Range("J14").Value = "ETHIAS,AME-LIFE,RESA,PROVINCE LIEGE"

# Widen column J (company breakdown) to fit the longer text, and move the
# active selection to below the used range, as in the saved workbook.
$ws.Columns.Item(10).ColumnWidth = 62
$ws.Range("J18").Select()
